# Update the "Förändrad" (Changed) date column (column C) for rows 2-28
# from serial date 45529 (2024-08-25) to 45530 (2024-08-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45529) {
        $cell.Value2 = 45530
    }
}
